# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) on affected leve rows across multiple sheets
# with freshly pulled values. A few rows also drop their now-unused
# NQ/HQ profit cell (ClearContents) where the source feed stopped
# reporting that split.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 360.31818
$ws.Range("I33").Value = 296.5
$ws.Range("K33").Value = 296.5
$ws.Range("M33").Value = -67.5

$ws.Range("H112").Value = 2664.1765
$ws.Range("J112").Value = 2664.1765
$ws.Range("L112").Value = 7992.529500000001
$ws.Range("N112").Value = -10208.5295

$ws.Range("H131").Value = 22653.8
$ws.Range("I131").Value = 22653.8
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 67961.39999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -62921.39999999999

$ws.Range("H137").Value = 10224.515
$ws.Range("I137").Value = 2933.3333
$ws.Range("J137").Value = 12748.385
$ws.Range("K137").Value = 8799.999899999999
$ws.Range("L137").Value = 38245.155
$ws.Range("M137").Value = -6249.999899999999
$ws.Range("N137").Value = -43345.155

$ws.Range("H138").Value = 5900.6924
$ws.Range("I138").Value = 7501.4443
$ws.Range("J138").Value = 5420.467
$ws.Range("K138").Value = 22504.3329
$ws.Range("L138").Value = 16261.401
$ws.Range("M138").Value = -17364.3329
$ws.Range("N138").Value = -26541.401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1345273
$ws.Range("I32").Value = 598921.75
$ws.Range("K32").Value = 598921.75
$ws.Range("M32").Value = -598634.75

$ws.Range("H61").Value = 3618.3635
$ws.Range("I61").Value = 2723.375
$ws.Range("J61").Value = 6005
$ws.Range("K61").Value = 2723.375
$ws.Range("L61").Value = 6005
$ws.Range("M61").Value = -2511.375
$ws.Range("N61").Value = -6429

$ws.Range("H88").Value = 4451.1665
$ws.Range("I88").Value = 2380.4
$ws.Range("J88").Value = 5930.2856
$ws.Range("K88").Value = 2380.4
$ws.Range("L88").Value = 5930.2856
$ws.Range("M88").Value = -1974.4
$ws.Range("N88").Value = -6742.2856

$ws.Range("H91").Value = 4451.1665
$ws.Range("I91").Value = 2380.4
$ws.Range("J91").Value = 5930.2856
$ws.Range("K91").Value = 2380.4
$ws.Range("L91").Value = 5930.2856
$ws.Range("M91").Value = -976.4000000000001
$ws.Range("N91").Value = -8738.285599999999

$ws.Range("H136").Value = 3618.3635
$ws.Range("I136").Value = 2723.375
$ws.Range("J136").Value = 6005
$ws.Range("K136").Value = 8170.125
$ws.Range("L136").Value = 18015
$ws.Range("M136").Value = -5620.125
$ws.Range("N136").Value = -23115

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2455049.5
$ws.Range("I31").Value = 3813.5
$ws.Range("J31").Value = 2911093.2
$ws.Range("K31").Value = 3813.5
$ws.Range("L31").Value = 2911093.2
$ws.Range("M31").Value = -3518.5
$ws.Range("N31").Value = -2911683.2

$ws.Range("H34").Value = 2455049.5
$ws.Range("I34").Value = 3813.5
$ws.Range("J34").Value = 2911093.2
$ws.Range("K34").Value = 3813.5
$ws.Range("L34").Value = 2911093.2
$ws.Range("M34").Value = -3611.5
$ws.Range("N34").Value = -2911497.2

$ws.Range("H134").Value = 3769.6897
$ws.Range("I134").Value = 4264.5264
$ws.Range("J134").Value = 2829.5
$ws.Range("K134").Value = 12793.5792
$ws.Range("L134").Value = 8488.5
$ws.Range("M134").Value = -10258.5792
$ws.Range("N134").Value = -13558.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1065.909
$ws.Range("I5").Value = 845.2
$ws.Range("J5").Value = 1249.8334
$ws.Range("K5").Value = 2535.6
$ws.Range("L5").Value = 3749.5002
$ws.Range("M5").Value = -2423.6
$ws.Range("N5").Value = -3973.5002

$ws.Range("H34").Value = 312
$ws.Range("I34").Value = 312
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 936
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -852

$ws.Range("H68").Value = 1567541.5
$ws.Range("J68").Value = 2005204.1
$ws.Range("L68").Value = 6015612.300000001
$ws.Range("N68").Value = -6017234.300000001

$ws.Range("H71").Value = 1567541.5
$ws.Range("J71").Value = 2005204.1
$ws.Range("L71").Value = 18046836.9
$ws.Range("N71").Value = -18054948.9

$ws.Range("H81").Value = 50
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 50
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = 150
$ws.Range("N81").Value = -2396

$ws.Range("H82").Value = 11601
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 11601
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").Value = 34803
$ws.Range("N82").Value = -35615

$ws.Range("H84").Value = 50
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 50
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = 450
$ws.Range("N84").Value = -11682

$ws.Range("H85").Value = 11601
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 11601
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("M85").Value = 34803
$ws.Range("N85").Value = -37611

$ws.Range("H104").Value = 5000
$ws.Range("I104").Value = 5000
$ws.Range("K104").Value = 15000
$ws.Range("M104").Value = -12379

$ws.Range("H107").Value = 3312.75
$ws.Range("J107").Value = 3643.2856
$ws.Range("L107").Value = 10929.8568
$ws.Range("N107").Value = -14769.8568

$ws.Range("H132").Value = 967.5
$ws.Range("I132").Value = 881
$ws.Range("J132").Value = 984.8
$ws.Range("K132").Value = 7929
$ws.Range("L132").Value = 8863.199999999999
$ws.Range("M132").Value = -5399
$ws.Range("N132").Value = -13923.2

$ws.Range("H135").Value = 1065.909
$ws.Range("I135").Value = 845.2
$ws.Range("J135").Value = 1249.8334
$ws.Range("K135").Value = 7606.8
$ws.Range("L135").Value = 11248.5006
$ws.Range("M135").Value = -5071.8
$ws.Range("N135").Value = -16318.5006

$ws.Range("H139").Value = 9744.277
$ws.Range("I139").Value = 13185
$ws.Range("J139").Value = 5443.375
$ws.Range("K139").Value = 39555
$ws.Range("L139").Value = 16330.125
$ws.Range("M139").Value = -34415
$ws.Range("N139").Value = -26610.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111113064
$ws.Range("I80").Value = 166667950
$ws.Range("J80").Value = 3300
$ws.Range("K80").Value = 166667950
$ws.Range("L80").Value = 3300
$ws.Range("M80").Value = -166666952
$ws.Range("N80").Value = -5296

$ws.Range("H83").Value = 111113064
$ws.Range("I83").Value = 166667950
$ws.Range("J83").Value = 3300
$ws.Range("K83").Value = 833339750
$ws.Range("L83").Value = 16500
$ws.Range("M83").Value = -833334758
$ws.Range("N83").Value = -26484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2520
$ws.Range("I122").Value = 2644.6
$ws.Range("J122").Value = 2208.5
$ws.Range("K122").Value = 7933.799999999999
$ws.Range("L122").Value = 6625.5
$ws.Range("M122").Value = -5483.799999999999
$ws.Range("N122").Value = -11525.5

$ws.Range("H135").Value = 86000
$ws.Range("J135").Value = 86000
$ws.Range("L135").Value = 86000
$ws.Range("N135").Value = -96140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49211.5
$ws.Range("J46").Value = 49211.5
$ws.Range("L46").Value = 49211.5
$ws.Range("N46").Value = -49673.5

$ws.Range("H122").Value = 50002800
$ws.Range("I122").Value = 2999
$ws.Range("K122").Value = 8997
$ws.Range("M122").Value = -6547

$ws.Range("H134").Value = 49211.5
$ws.Range("J134").Value = 49211.5
$ws.Range("L134").Value = 147634.5
$ws.Range("N134").Value = -152704.5
